# Apply cryptos list update (prices & volume %) per commit: "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'61.759.48"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3
$ws.Range("D3").Value = "'2.887.10"
$ws.Range("E3").Value = "  -2.07%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'566.77"
$ws.Range("E5").Value = "  -3.71%  "

# Row 6
$ws.Range("D6").Value = "'142.85"
$ws.Range("E6").Value = "  -2.53%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.99%  "

# Row 9
$ws.Range("D9").Value = "'2.884.60"
$ws.Range("E9").Value = "  -2.14%  "

# Row 10
$ws.Range("D10").Value = "'6.98"
$ws.Range("E10").Value = "  +0.17%  "

# Row 11
$ws.Range("D11").Value = "'0.146"
$ws.Range("E11").Value = "  -1.55%  "

# Row 12
$ws.Range("D12").Value = "'0.429"
$ws.Range("E12").Value = "  -1.31%  "

# Row 13
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  -0.46%  "

# Row 14
$ws.Range("D14").Value = "'31.89"
$ws.Range("E14").Value = "  -0.75%  "

# Row 15
$ws.Range("E15").Value = "  +0.02%  "

# Row 16
$ws.Range("D16").Value = "'3.368.11"
$ws.Range("E16").Value = "  -2.10%  "

# Row 17
$ws.Range("D17").Value = "'61.727.10"
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'2.897.45"
$ws.Range("E18").Value = "  -1.79%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.52"
$ws.Range("E19").Value = "  -1.81%  "

# Row 20
$ws.Range("D20").Value = "'428.53"
$ws.Range("E20").Value = "  -1.27%  "

# Row 21
$ws.Range("D21").Value = "'12.96"
$ws.Range("E21").Value = "  -3.60%  "

# Row 22
$ws.Range("D22").Value = "'0.652"
$ws.Range("E22").Value = "  -1.10%  "

# Row 23
$ws.Range("D23").Value = "'6.85"
$ws.Range("E23").Value = "  -1.33%  "

# Row 24
$ws.Range("D24").Value = "'78.83"
$ws.Range("E24").Value = "  -1.65%  "

# Row 25
$ws.Range("D25").Value = "'12.00"

# Row 26
$ws.Range("D26").Value = "'10.03"
$ws.Range("E26").Value = "  -9.77%  "

# Row 27
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("E28").Value = "  -3.53%  "

# Row 29
$ws.Range("E29").Value = "  +9.85%  "

# Row 30
$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  -2.67%  "

# Row 31
$ws.Range("E31").Value = "  -3.18%  "

# Row 32
$ws.Range("D32").Value = "'2.02"
$ws.Range("E32").Value = "  -6.68%  "

# Row 33
$ws.Range("E33").Value = "  -0.11%  "

# Row 34
$ws.Range("E34").Value = "  -1.16%  "

# Row 35
$ws.Range("D35").Value = "'25.56"
$ws.Range("E35").Value = "  -1.99%  "

# Row 36
$ws.Range("E36").Value = "  -4.43%  "

# Row 37
$ws.Range("E37").Value = "  -3.28%  "

# Row 38
$ws.Range("D38").Value = "'48.85"
$ws.Range("E38").Value = "  -1.59%  "

# Row 39
$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "  -6.40%  "

# Row 40
$ws.Range("E40").Value = "  -4.50%  "

# Row 41
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
$ws.Range("D42").Value = "'8.13"
$ws.Range("E42").Value = "  -2.38%  "

# Row 43
$ws.Range("D43").Value = "'39.76"
$ws.Range("E43").Value = "  +1.76%  "

# Row 44
$ws.Range("E44").Value = "  -2.14%  "

# Row 45
$ws.Range("D45").Value = "'2.689.48"
$ws.Range("E45").Value = "  +0.42%  "

# Row 46
$ws.Range("E46").Value = "  +0.34%  "

# Row 47
$ws.Range("E47").Value = "  -2.58%  "

# Row 48
$ws.Range("D48").Value = "'344.72"
$ws.Range("E48").Value = "  -2.11%  "

# Row 51
$ws.Range("D51").Value = "'21.55"
$ws.Range("E51").Value = "  -4.04%  "
